$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: adm_url_source -> URL (as hyperlink, styled like existing hyperlink in B2)
$ws.Range("A8").Value = "adm_url_source"
$ws.Range("B8").Value = "https://geoportal.dane.gov.co/descargas/veredas/CRVeredas_2017.zip"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://geoportal.dane.gov.co/descargas/veredas/CRVeredas_2017.zip") | Out-Null
$ws.Range("B8").Style = "Hyperlink"

# Row 9: plot_character_file -> ;
$ws.Range("A9").Value = "plot_character_file"
$ws.Range("B9").Value = ";"

$wb.Save()
